# Updated cryptos list on Sat Oct 14 15:28:06 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row. Both columns hold text in the original workbook (prices use
# a "thousands.thousands.decimal" style grouping and percentages keep
# leading/trailing padding spaces), so for any new Price value that would
# otherwise be auto-parsed by Excel as a genuine number we force the cell to
# text (NumberFormat "@") before writing it and then reset the cell style
# back to "Normal" so no stray number-format style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.944.62"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.554.58"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.775.07"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.546.32"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "26.931.20"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("E18").Value = "  +2.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("D34").Value = "1.417.90"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  +13.92%  "
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.528"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "1.688.43"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0521"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").Value = "0.0₇0997"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0960"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.76%  "
